$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-empty "Quantity" column (C) for each part row ---
$ws.Range("C2").Value = 6
$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 2

# --- Highlight rows: green for normal parts, yellow for the diode (D1) row ---
$ws.Range("A2:C2").Interior.Color = 0x50B000
$ws.Range("A3:C3").Interior.Color = 0x50B000
$ws.Range("A4:C4").Interior.Color = 0x50B000
$ws.Range("A5:C5").Interior.Color = 0x50B000
$ws.Range("A7:C7").Interior.Color = 0x50B000
$ws.Range("A6:C6").Interior.Color = 0xFFFF

$ws.Range("A2:C7").HorizontalAlignment = -4108
$ws.Range("A2:C7").VerticalAlignment = -4108

# --- MOUSER part number + URL for the diode (D1) ---
$ws.Range("G6").Value = "494-SMBJ4729AE3/TR13 "
$ws.Range("G6").Style = "Normal"
$ws.Hyperlinks.Add($ws.Range("H6"), "https://hr.mouser.com/ProductDetail/Microsemi/SMBJ4729Ae3-TR13/?qs=%2fha2pyFadugozOLVQX8oU%252bLu9j3Llt5qc4Sq4QXvUCgpmf6R2dENSw%3d%3d")

# --- Column width for the now-wider MOUSER part no. column ---
$ws.Columns.Item(7).ColumnWidth = 21

# --- Selection moved to G7 (as left by the editor) ---
$ws.Range("G7").Select()
